# Apply the "Map032 scene update" edit:
#  - Add a new value "Rewrite     -   Lily" in cell C2 (next to "Lily" row, A2)
#  - Add a new value "Sina" in cell C3 (next to "Sina" row, A3)
# The existing column A values are left untouched; only the new column C
# cells are introduced, which also extends the used range to A1:C46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Rewrite     -   Lily"
$ws.Range("C3").Value = "Sina"
